$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (was "GossF-HW50.xpc" -> "GossF")
$ws.Name = "GossF"

# Append a new row (16) that repeats the "HexGrid-60degTilt5degRes" scheme
# row with index 14, mirroring the pattern of rows 2-15.
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16:M16").Value = 1

# Copy the formatting from row 15 (A15 label cell, C15:M15 data cells)
# onto the new row so the new row matches the existing style pattern.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C15:M15").Copy()
$ws.Range("C16:M16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
